# Weekly fruit/vegetable price update: a new week's record is inserted at
# row 37 (Poroto verde, "Sin especificar" variety, Región del Maule origin),
# pushing the previously-existing rows 37-68 down to 38-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 37; Excel shifts rows 37:68 down to 38:69
# and grows the used range to A1:R69 automatically.
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with this week's record.
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 44603
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 60
$ws.Range("K37").Value = 35000
$ws.Range("L37").Value = 36000
$ws.Range("M37").Value = 35500
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 1420
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
